$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.370.39'
$ws.Range('E2').Value = '  +9.21%  '
$ws.Range('D3').Value = '1.676.96'
$ws.Range('E3').Value = '  +4.94%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = "'306.65"
$ws.Range('E5').Value = '  +6.21%  '
$ws.Range('D6').Value = "'0.9976"
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').Value = "'48.09"
$ws.Range('E9').Value = '  +12.67%  '
$ws.Range('E10').Value = '  +3.55%  '
$ws.Range('D11').Value = "'0.07257"
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = "'20.32"
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').Value = "'6.097"
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').Value = "'6.742"
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '1.677.79'
$ws.Range('E16').Value = '  +5.02%  '
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('D18').Value = "'0.9980"
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').Value = "'0.06719"
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').Value = "'81.05"
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').Value = "'16.42"
$ws.Range('E21').Value = '  +1.87%  '
$ws.Range('D22').Value = "'6.093"
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').Value = "'11.93"
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').Value = '24.315.35'
$ws.Range('E24').Value = '  +8.93%  '
$ws.Range('D25').Value = "'2.431"
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').Value = "'3.364"
$ws.Range('E26').Value = '  -11.42%  '
$ws.Range('D27').Value = "'2.656"
$ws.Range('E27').Value = '  +6.81%  '
$ws.Range('D28').Value = "'152.21"
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').Value = "'19.56"
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '1.862.12'
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').Value = "'127.10"
$ws.Range('E31').Value = '  +5.52%  '
$ws.Range('E32').Value = '  +5.05%  '
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').Value = "'0.9659"
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('D35').Value = "'1.742"
$ws.Range('E35').Value = '  +8.04%  '
$ws.Range('D36').Value = "'0.08491"
$ws.Range('E36').Value = '  +2.68%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = "'9.038"
$ws.Range('E37').Value = '  +4.94%  '
$ws.Range('D38').Value = "'0.06462"
$ws.Range('E38').Value = '  +5.19%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = "'12.29"
$ws.Range('E39').Value = '  +4.68%  '
$ws.Range('D40').Value = "'5.330"
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = "'0.02327"
$ws.Range('E41').Value = '  +5.29%  '
$ws.Range('D42').Value = "'1.260"
$ws.Range('E42').Value = '  +2.22%  '
$ws.Range('D43').Value = "'0.2106"
$ws.Range('E43').Value = '  +4.01%  '
$ws.Range('D44').Value = "'0.6167"
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('D45').Value = "'0.9975"
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').Value = "'3.775"
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'13.06"
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = "'0.5943"
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('D49').Value = "'127.07"
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = "'2.022"
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('D51').Value = "'0.07217"
$ws.Range('E51').Value = '  +5.84%  '
